$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "68.812.44"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -3.64%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.493.29"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -4.20%  "

$ws.Range("E4").Value = "  -0.15%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "577.82"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "176.79"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.76%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.619"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.37%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.487.53"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -4.14%  "

$ws.Range("E9").Value = "  -0.08%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.188"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -6.86%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.65"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +7.78%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.603"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "47.30"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -4.73%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000277"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "687.50"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.01%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "8.90"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.16%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.046.00"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.47%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "68.747.22"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.01%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.486.48"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.16%  "

$ws.Range("E20").Value = "  -1.71%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.52"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -3.67%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "11.17"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.91%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.904"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -3.83%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "16.33"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -8.56%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "97.92"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.13%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.84"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.53%  "

$ws.Range("E27").Value = "  +0.14%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.66"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -6.11%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.41"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -7.43%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "33.03"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -6.31%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "8.76"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -4.41%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -7.76%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.31"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.62%  "

$ws.Range("E34").Value = "  -6.37%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "565.44"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.67"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -12.87%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "10.92"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.47%  "

$ws.Range("E38").Value = "  -3.06%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "56.73"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -4.65%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0440"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.76%  "

$ws.Range("E42").Value = "  -4.06%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.336"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.52%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.421.99"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -8.59%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "33.47"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.78%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0$([char]0x2083)0703"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -7.64%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.94"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.31%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.61"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.42%  "

$ws.Range("E49").Value = "  -0.37%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "134.26"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("E51").Value = "  -0.94%  "
